$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1842, 81.599999999999994),
    @(2097.1999999999998, 81.400000000000006),
    @(2352.4, 81.2),
    @(2607.6, 81),
    @(2862.8, 80.8),
    @(3118, 80.599999999999994),
    @(3373.2, 80.400000000000006),
    @(3628.4, 80.2),
    @(3883.6, 80),
    @(4138.8, 79.8),
    @(4394, 79.599999999999994),
    @(4649.2, 79.400000000000006),
    @(4904.3999999999996, 79.2),
    @(5159.6000000000004, 79),
    @(5414.8, 78.8),
    @(5670, 78.599999999999994),
    @(5925.2, 78.399999999999906),
    @(6180.4, 78.199999999999903),
    @(6435.6, 77.999999999999901),
    @(6690.8, 77.799999999999898),
    @(6946, 77.599999999999895),
    @(7201.2, 77.399999999999906),
    @(7456.4, 77.199999999999903),
    @(7711.6, 76.999999999999901),
    @(7966.8, 76.799999999999898),
    @(8222, 76.599999999999895),
    @(8477.2000000000007, 76.399999999999906),
    @(8732.4, 76.199999999999903),
    @(8987.6, 75.999999999999901),
    @(9242.7999999999993, 75.799999999999898),
    @(9498, 75.599999999999895),
    @(9753.2000000000007, 75.399999999999906),
    @(10008.4, 75.199999999999903),
    @(10263.6, 74.999999999999901),
    @(10518.8, 74.799999999999898),
    @(10774, 74.599999999999895),
    @(11029.2, 74.399999999999906),
    @(11284.4, 74.199999999999903),
    @(11539.6, 73.999999999999901),
    @(11794.8, 73.799999999999898),
    @(12050, 73.599999999999895),
    @(12305.2, 73.399999999999906),
    @(12560.4, 73.199999999999903),
    @(12815.6, 72.999999999999901),
    @(13070.8, 72.799999999999898),
    @(13326, 72.599999999999895),
    @(13581.2, 72.399999999999906),
    @(13836.4, 72.199999999999903),
    @(14091.6, 71.999999999999901),
    @(14346.8, 71.799999999999898),
    @(14602, 71.599999999999895),
    @(14857.2, 71.399999999999807),
    @(15112.4, 71.199999999999804),
    @(15367.6, 70.999999999999801),
    @(15622.8, 70.799999999999798),
    @(15878, 70.599999999999795),
    @(16133.2, 70.399999999999807),
    @(16388.400000000001, 70.199999999999804),
    @(16643.599999999999, 69.999999999999801)
)

$startRow = 17
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

$ws.Range("E74").Select()
